$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "28.883.26"
$ws.Range("E2").Value = "  +1.68%  "

# Row 3
$ws.Range("D3").Value = "1.900.29"
$ws.Range("E3").Value = "  +1.52%  "

# Row 4
$ws.Range("D4").Formula = "'1.002"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.16%  "

# Row 5
$ws.Range("D5").Formula = "'332.50"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -1.75%  "

# Row 6
$ws.Range("D6").Formula = "'1.001"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.06%  "

# Row 7
$ws.Range("D7").Formula = "'0.4628"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -1.30%  "

# Row 8
$ws.Range("D8").Formula = "'0.4071"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +2.79%  "

# Row 9
$ws.Range("D9").Formula = "'47.97"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +0.28%  "

# Row 10
$ws.Range("D10").Formula = "'0.07997"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -0.42%  "

# Row 11
$ws.Range("D11").Formula = "'1.004"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.08%  "

# Row 12
$ws.Range("D12").Formula = "'21.71"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -0.87%  "

# Row 13
$ws.Range("D13").Value = "1.906.51"
$ws.Range("E13").Value = "  +1.97%  "

# Row 14
$ws.Range("D14").Formula = "'5.938"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -0.91%  "

# Row 15
$ws.Range("D15").Formula = "'7.092"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -2.11%  "

# Row 16
$ws.Range("D16").Formula = "'89.14"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -2.11%  "

# Row 17
$ws.Range("D17").Formula = "'1.001"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -0.14%  "

# Row 18
$ws.Range("D18").Formula = "'0.00001033"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -0.67%  "

# Row 19
$ws.Range("D19").Formula = "'0.06560"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.79%  "

# Row 20
$ws.Range("D20").Formula = "'17.48"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -0.59%  "

# Row 21
$ws.Range("D21").Formula = "'1.004"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.38%  "

# Row 22
$ws.Range("D22").Value = "28.900.19"
$ws.Range("E22").Value = "  +1.72%  "

# Row 23
$ws.Range("D23").Formula = "'5.469"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.29%  "

# Row 24
$ws.Range("D24").Formula = "'11.08"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.32%  "

# Row 25
$ws.Range("E25").Value = "  -1.44%  "

# Row 26
$ws.Range("D26").Value = "2.132.85"
$ws.Range("E26").Value = "  +1.86%  "

# Row 27
$ws.Range("D27").Formula = "'157.34"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -1.86%  "

# Row 28
$ws.Range("D28").Formula = "'19.74"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -0.08%  "

# Row 29
$ws.Range("D29").Formula = "'2.100"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -1.74%  "

# Row 30
$ws.Range("D30").Formula = "'5.384"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -2.30%  "

# Row 31
$ws.Range("D31").Formula = "'118.92"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -0.96%  "

# Row 32
$ws.Range("D32").Formula = "'0.9823"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +1.16%  "

# Row 33
$ws.Range("D33").Formula = "'0.09382"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -1.13%  "

# Row 34
$ws.Range("D34").Formula = "'1.414"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +2.50%  "

# Row 35
$ws.Range("D35").Formula = "'3.596"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +0.73%  "

# Row 36
$ws.Range("D36").Formula = "'5.302"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -0.95%  "

# Row 37
$ws.Range("D37").Formula = "'0.06080"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -0.50%  "

# Row 38
$ws.Range("D38").Formula = "'0.02227"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -1.28%  "

# Row 39
$ws.Range("D39").Formula = "'8.405"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +0.35%  "

# Row 40
$ws.Range("D40").Formula = "'1.169"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -1.27%  "

# Row 41
$ws.Range("D41").Formula = "'0.5828"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -2.08%  "

# Row 42
$ws.Range("D42").Formula = "'1.001"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +0.09%  "

# Row 43
$ws.Range("D43").Formula = "'10.16"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -1.76%  "

# Row 44
$ws.Range("D44").Formula = "'0.1826"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -2.56%  "

# Row 45
$ws.Range("D45").Formula = "'1.260"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -2.66%  "

# Row 46
$ws.Range("D46").Formula = "'2.358"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +15.18%  "

# Row 47
$ws.Range("D47").Formula = "'12.17"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -0.26%  "

# Row 48
$ws.Range("D48").Formula = "'0.5495"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -1.45%  "

# Row 49
$ws.Range("D49").Formula = "'1.904"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -2.85%  "

# Row 50
$ws.Range("D50").Formula = "'0.07027"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +2.37%  "

# Row 51
$ws.Range("D51").Formula = "'47.10"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +21.90%  "
